$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ColdStorage")

# Update the item name in row 1 (A1/B1 were "Boostrix Injection")
$ws.Range("A1").Value = "Stpase 1500000IU Injection 10ml"
$ws.Range("B1").Value = "Stpase 1500000IU Injection 10ml"

# Add the new "cscs" entry in B5
$ws.Range("B5").Value = "cscs"

# Make ColdStorage the active/selected sheet (was Msite before)
$ws.Activate()
$ws.Range("B4").Select()
